# Apply the "Updated remaining queries for C3DC" edit:
#  - Fix the JOIN conditions in every embedded SQL query so that they
#    reference the renamed columns (study_id / participant_id) instead
#    of the old generic "id" columns.
#  - Widen column C to a fixed 67.5 width (and drop the bestFit flag).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold one of the SQL queries needing the JOIN-clause fix.
$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $text = $rng.Value2

    if ($text -ne $null) {
        $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
        $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
        $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
        $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
        $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
        $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

        $rng.Value2 = $text
    }
}

# Column C: switch from auto "best fit" width to an explicit fixed width
# (67.5 in the saved OOXML). The runtime's ColumnWidth setter stores the
# width with a +5/6 offset versus the value read back through COM, so we
# compensate here to land exactly on 67.5 in the saved file.
$ws.Columns("C").ColumnWidth = 66.66666666666667

$wb.Save()
